# Converts a hex RRGGBB string into the little-endian integer that the
# PowerPoint object model expects for a .RGB property (same encoding as
# VBA's RGB() function: R + G*256 + B*65536).
function Hex2RGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-theme the deck: the design was switched from the "Integral" /
#    "Red Violet" colour scheme back to the stock Office colour scheme.
#    The 12 theme colours are a per-package resource shared by every
#    slide (they all hang off the one slide master), so touching them
#    from any single slide updates them everywhere.
# ---------------------------------------------------------------------
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$firstSlide = $p.Slides.Item(1)
$themeColors = $firstSlide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Hex2RGB $officeColors[$i - 1]
}

# ---------------------------------------------------------------------
# 2) Re-style the three tables that were bound to the old theme's table
#    style ("Table_0") so they pick up the built-in table style instead.
# ---------------------------------------------------------------------
$newTableStyleId = "{1B770D39-F8F5-4B8A-BEE8-3404291CDD1E}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
